$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A1").Value = "Коэффициент скорости обучения"
$ws.Range("B1").Value = 0.1

$ws.Range("A2").Value = "Количество входов нейронной сети"
$ws.Range("B2").Value = 4

$ws.Range("A3").Value = "Размерность выходного слоя"
$ws.Range("B3").Value = 2

$ws.Range("G8").Select()
